$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Gắn giao diện BNS và BDT": the notification list now only carries the
# Nhân sự (HR) and Đào tạo (Training) department notices, so the two
# in-between rows (id=2 "Thông báo lịch làm việc 1-1" and id=3 "Thông báo
# đến dọn xưởng") are removed. The remaining "Thông báo tuyển dụng" /
# "Ban Đào tạo" row shifts up to become row 3.
$ws.Rows("3:4").Delete()

# Touch the header/footer settings so the worksheet carries an explicit
# (empty) headerFooter section, matching the resaved sheet.
$ws.PageSetup.CenterHeader = ""
